$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for id_sector 4 and 5 (original rows 28:29) first (from bottom up)
# to avoid shifting the indices of the rows we still need to delete.
$ws.Rows("28:29").Delete() | Out-Null

# Delete rows for id_sector 1 and 2 (original rows 2:11)
$ws.Rows("2:11").Delete() | Out-Null

# Update the selection to match the target state
$ws.Range("D24").Select() | Out-Null
